$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.868900000000005
$ws.Range("D9").Value = -8.117699999999999
$ws.Range("A11").Value = -21.9686
$ws.Range("C11").Value = -12.14820000000001
$ws.Range("A12").Value = -21.36269999999999
$ws.Range("D13").Value = -8.968399999999992
$ws.Range("D14").Value = -8.042400000000001
$ws.Range("A15").Value = -21.6728
$ws.Range("D19").Value = -7.972499999999999
$ws.Range("D21").Value = -8.702399999999994
$ws.Range("D22").Value = -7.958500000000003
$ws.Range("C23").Value = -12.09539999999999
$ws.Range("D24").Value = -7.332499999999999
$ws.Range("D26").Value = -7.721299999999996
$ws.Range("A27").Value = -21.8142
$ws.Range("A28").Value = -21.98520000000001
$ws.Range("C28").Value = -13.58229999999999
$ws.Range("A31").Value = -21.4925
$ws.Range("A32").Value = -21.7298
$ws.Range("C32").Value = -12.8039
$ws.Range("C34").Value = -10.74840000000002
$ws.Range("A36").Value = -19.86729999999998
$ws.Range("C36").Value = -12.68890000000001
$ws.Range("C37").Value = -13.80669999999999
$ws.Range("A38").Value = -20.05539999999997
$ws.Range("D38").Value = -8.88709999999999
$ws.Range("D41").Value = -8.304500000000001
$ws.Range("C42").Value = -13.1076
$ws.Range("A46").Value = -21.93060000000001
$ws.Range("C49").Value = -13.97100000000001
$ws.Range("D52").Value = -8.083499999999997
$ws.Range("A54").Value = -21.88750000000001
$ws.Range("C54").Value = -13.54029999999999
$ws.Range("A55").Value = -21.99370000000002
$ws.Range("A56").Value = -21.9653
$ws.Range("D56").Value = -8.887499999999999
$ws.Range("A67").Value = -21.55039999999996
$ws.Range("A69").Value = -21.65739999999997
$ws.Range("D71").Value = -7.131599999999999
$ws.Range("A72").Value = -21.78379999999999
$ws.Range("D72").Value = -7.381400000000001
$ws.Range("A73").Value = -20.11880000000002
$ws.Range("C78").Value = -12.37280000000001
$ws.Range("D78").Value = -8.848699999999994
$ws.Range("C80").Value = -12.02
$ws.Range("A83").Value = -21.7106
$ws.Range("D83").Value = -8.095099999999999
$ws.Range("D85").Value = -9.019500000000004
$ws.Range("A86").Value = -21.98860000000002
$ws.Range("D86").Value = -8.352800000000002
$ws.Range("D90").Value = -6.625999999999994
$ws.Range("A91").Value = -20.51169999999998
$ws.Range("A93").Value = -21.55060000000002
$ws.Range("D96").Value = -7.960099999999994
$ws.Range("C97").Value = -10.5383
$ws.Range("A99").Value = -21.8288
$ws.Range("C99").Value = -11.8664
$ws.Range("C100").Value = -12.1922
$ws.Range("C101").Value = -12.5745
$ws.Range("D103").Value = -8.636699999999996
$ws.Range("A104").Value = -21.35039999999999
$ws.Range("A105").Value = -20.11729999999998
